$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161, shifting existing rows 161-223 down to 162-224.
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new data record.
$ws.Range("A161").Value = 9
$ws.Range("B161").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C161").Value = "Metropolitana"
$ws.Range("D161").Value = 44726
$ws.Range("E161").Value = 13
$ws.Range("F161").Value = 100112026
$ws.Range("G161").Value = "Haba"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 39
$ws.Range("K161").Value = 21000
$ws.Range("L161").Value = 22000
$ws.Range("M161").Value = 21487
$ws.Range("N161").Value = "$/saco 25 kilos"
$ws.Range("O161").Value = "Provincia del Elquí"
$ws.Range("P161").Value = 859
$ws.Range("Q161").Value = 25
$ws.Range("R161").Value = "Hortaliza"
